$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay stored as text (matching original inlineStr format)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '49.915.06'
$ws.Cells.Item(2, 5).Value = '  +3.83%  '

$ws.Cells.Item(3, 4).Value = '2.658.69'
$ws.Cells.Item(3, 5).Value = '  +6.31%  '

$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).Value = '327.97'
$ws.Cells.Item(5, 5).Value = '  +2.23%  '

$ws.Cells.Item(6, 4).Value = '111.37'
$ws.Cells.Item(6, 5).Value = '  +3.72%  '

$ws.Cells.Item(7, 4).Value = '0.531'
$ws.Cells.Item(7, 5).Value = '  +1.19%  '

$ws.Cells.Item(8, 5).Value = '  +0.03%  '

$ws.Cells.Item(9, 5).Value = '  +3.85%  '

$ws.Cells.Item(10, 4).Value = '40.79'
$ws.Cells.Item(10, 5).Value = '  +2.97%  '

$ws.Cells.Item(11, 4).Value = '20.65'
$ws.Cells.Item(11, 5).Value = '  +2.72%  '

$ws.Cells.Item(12, 5).Value = '  +1.46%  '

$ws.Cells.Item(13, 5).Value = '  +0.81%  '

$ws.Cells.Item(14, 4).Value = '7.32'
$ws.Cells.Item(14, 5).Value = '  +3.16%  '

$ws.Cells.Item(15, 4).Value = '3.075.97'
$ws.Cells.Item(15, 5).Value = '  +6.22%  '

$ws.Cells.Item(16, 4).Value = '2.659.04'
$ws.Cells.Item(16, 5).Value = '  +6.03%  '

$ws.Cells.Item(17, 5).Value = '  +6.15%  '

$ws.Cells.Item(18, 4).Value = '49.936.99'
$ws.Cells.Item(18, 5).Value = '  +4.08%  '

$ws.Cells.Item(19, 5).Value = '  +3.00%  '

$ws.Cells.Item(20, 4).Value = '2.99'
$ws.Cells.Item(20, 5).Value = '  +8.33%  '

$ws.Cells.Item(21, 4).Value = '6.86'

$ws.Cells.Item(22, 4).Value = '0.0₃0964'
$ws.Cells.Item(22, 5).Value = '  +2.64%  '

$ws.Cells.Item(23, 2).Value = 'Litecoin'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(23, 4).Value = '73.26'
$ws.Cells.Item(23, 5).Value = '  +2.55%  '

$ws.Cells.Item(24, 2).Value = 'BitcoinCash'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(24, 4).Value = '282.11'
$ws.Cells.Item(24, 5).Value = '  +1.93%  '

$ws.Cells.Item(25, 5).Value = '  +2.63%  '

$ws.Cells.Item(26, 4).Value = '27.08'
$ws.Cells.Item(26, 5).Value = '  +4.56%  '

$ws.Cells.Item(27, 5).Value = '  -0.01%  '

$ws.Cells.Item(28, 5).Value = '  +6.96%  '

$ws.Cells.Item(29, 5).Value = '  +2.87%  '

$ws.Cells.Item(30, 4).Value = '36.74'
$ws.Cells.Item(30, 5).Value = '  +4.43%  '

$ws.Cells.Item(31, 5).Value = '  +2.23%  '

$ws.Cells.Item(32, 4).Value = '49.75'
$ws.Cells.Item(32, 5).Value = '  +0.13%  '

$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).Value = '5.46'
$ws.Cells.Item(33, 5).Value = '  +3.08%  '

$ws.Cells.Item(34, 2).Value = 'Celestia'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(34, 4).Value = '19.61'
$ws.Cells.Item(34, 5).Value = '  +0.56%  '

$ws.Cells.Item(35, 5).Value = '  -0.07%  '

$ws.Cells.Item(36, 5).Value = '  +1.90%  '

$ws.Cells.Item(37, 5).Value = '  +7.28%  '

$ws.Cells.Item(38, 5).Value = '  +3.10%  '

$ws.Cells.Item(39, 5).Value = '  +8.49%  '

$ws.Cells.Item(40, 4).Value = '126.51'
$ws.Cells.Item(40, 5).Value = '  +4.70%  '

$ws.Cells.Item(41, 5).Value = '  +1.80%  '

$ws.Cells.Item(42, 4).Value = '22.73'
$ws.Cells.Item(42, 5).Value = '  +7.07%  '

$ws.Cells.Item(43, 5).Value = '  +0.95%  '

$ws.Cells.Item(44, 4).Value = '0.0314'
$ws.Cells.Item(44, 5).Value = '  +3.96%  '

$ws.Cells.Item(45, 5).Value = '  +7.31%  '

$ws.Cells.Item(46, 4).Value = '2.069.76'
$ws.Cells.Item(46, 5).Value = '  +2.34%  '

$ws.Cells.Item(47, 5).Value = '  +14.08%  '

$ws.Cells.Item(48, 4).Value = '1.99'
$ws.Cells.Item(48, 5).Value = '  +8.15%  '

$ws.Cells.Item(49, 4).Value = '9.09'
$ws.Cells.Item(49, 5).Value = '  +1.26%  '

$ws.Cells.Item(50, 4).Value = '5.43'
$ws.Cells.Item(50, 5).Value = '  +4.98%  '

$ws.Cells.Item(51, 4).Value = '81.85'
$ws.Cells.Item(51, 5).Value = '  +1.91%  '

# Restore default General format/style so only cell content differs from the original
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
